$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers
$ws.Range("K1").Value = "Payment Status"
$ws.Range("L1").Value = "Amount Paid"
$ws.Range("K1:L1").Style = $ws.Range("A1").Style

# Set empty strings for existing rows 2-10 in K and L columns
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 11).Value = ""
    $ws.Cells.Item($r, 12).Value = ""
}

# Add row 11: duplicate of pattern
$ws.Range("A11").Value = "2025-12-31"
$ws.Range("B11").Value = "Wednesday"
$ws.Range("C11").Value = "vairgwadi"
$ws.Range("D11").Value = "Vinayak"
$ws.Range("E11").Value = "GOLD Tea Powder"
$ws.Range("F11").Value = "Mix"
$ws.Range("G11").Value = "100gm"
$ws.Range("H11").Value = 60
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 60
$ws.Range("K11").Value = ""
$ws.Range("L11").Value = ""

# Add row 12
$ws.Range("A12").Value = "2025-12-31"
$ws.Range("B12").Value = "Wednesday"
$ws.Range("C12").Value = "vairgwadi"
$ws.Range("D12").Value = "Vinayak"
$ws.Range("E12").Value = "GOLD Tea Powder"
$ws.Range("F12").Value = "Mix"
$ws.Range("G12").Value = "100gm"
$ws.Range("H12").Value = 60
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 60
$ws.Range("K12").Value = "Half paid"
$ws.Range("L12").Value = 50
